# Update USB connector (CON1) part info: replace DX4R005JJ7R1500 (JAE Electronics,
# 670-2678-1-ND) with UJ2-MIBH-4-MSMT-TR (CUI Inc., 102-4008-1-ND)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row for designator CON1 in column A so we don't hardcode row numbers
$designatorCell = $ws.Columns.Item(1).Find("CON1")

if ($designatorCell -ne $null) {
    $row = $designatorCell.Row

    # Remember a cell in this row whose formatting (quote-prefixed text style)
    # matches what all the part-info cells should keep, so we can restore it
    # after writing new values (assigning .Value alone resets the cell style).
    $fmtSource = $ws.Cells.Item($row, 1)

    # Comment column (C) and Manufacturer Part Number column (E) both held the old MPN
    $ws.Cells.Item($row, 3).Value = "UJ2-MIBH-4-MSMT-TR"
    # Manufacturer column (D)
    $ws.Cells.Item($row, 4).Value = "CUI Inc."
    # Manufacturer Part Number column (E)
    $ws.Cells.Item($row, 5).Value = "UJ2-MIBH-4-MSMT-TR"
    # Digi-Key Part Number column (F)
    $ws.Cells.Item($row, 6).Value = "102-4008-1-ND"

    # Restore original cell formatting/style on the edited cells (C:F)
    $fmtSource.Copy()
    $targetRange = $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 6))
    $targetRange.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}
